$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5525.25
$ws.Range("J43").Value = 4050
$ws.Range("L43").Value = 4050
$ws.Range("N43").Value = -4188
$ws.Range("H74").Value = 6273.2
$ws.Range("I74").Value = 5091.5
$ws.Range("K74").Value = 5091.5
$ws.Range("M74").Value = -4155.5
$ws.Range("H77").Value = 6273.2
$ws.Range("I77").Value = 5091.5
$ws.Range("K77").Value = 25457.5
$ws.Range("M77").Value = -20777.5
$ws.Range("H100").Value = 1921
$ws.Range("J100").Value = 1959
$ws.Range("L100").Value = 1959
$ws.Range("N100").Value = -3041
$ws.Range("H115").Value = 2852
$ws.Range("I115").Value = 2852
$ws.Range("K115").Value = 8556
$ws.Range("M115").Value = -6989
$ws.Range("H116").Value = 4599.273
$ws.Range("I116").Value = 4766
$ws.Range("K116").Value = 4766
$ws.Range("M116").Value = -1324
$ws.Range("H138").Value = 2144.8333
$ws.Range("I138").Value = 1233
$ws.Range("J138").Value = 2796.1428
$ws.Range("K138").Value = 3699
$ws.Range("L138").Value = 8388.428400000001
$ws.Range("M138").Value = 1441
$ws.Range("N138").Value = -18668.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5793.234
$ws.Range("I32").Value = 5875.7173
$ws.Range("K32").Value = 5875.7173
$ws.Range("M32").Value = -5588.7173
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2715.3333
$ws.Range("J64").Value = 2949.4285
$ws.Range("L64").Value = 2949.4285
$ws.Range("N64").Value = -3399.4285
$ws.Range("H67").Value = 2715.3333
$ws.Range("J67").Value = 2949.4285
$ws.Range("L67").Value = 2949.4285
$ws.Range("N67").Value = -4509.4285
$ws.Range("H105").Value = 3381.25
$ws.Range("I105").Value = 2087.8572
$ws.Range("K105").Value = 2087.8572
$ws.Range("M105").Value = -340.8571999999999
$ws.Range("H134").Value = 2779.05
$ws.Range("I134").Value = 2643.5
$ws.Range("K134").Value = 7930.5
$ws.Range("M134").Value = -5395.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 33990
$ws.Range("I22").Value = 984.5
$ws.Range("J22").Value = 100001
$ws.Range("K22").Value = 984.5
$ws.Range("L22").Value = 100001
$ws.Range("M22").Value = -634.5
$ws.Range("N22").Value = -100701
$ws.Range("H31").Value = 1592.641
$ws.Range("I31").Value = 1626.2084
$ws.Range("K31").Value = 1626.2084
$ws.Range("M31").Value = -1331.2084
$ws.Range("H34").Value = 1592.641
$ws.Range("I34").Value = 1626.2084
$ws.Range("K34").Value = 1626.2084
$ws.Range("M34").Value = -1424.2084
$ws.Range("H58").Value = 2745.2222
$ws.Range("I58").Value = 2650.875
$ws.Range("K58").Value = 2650.875
$ws.Range("M58").Value = -2447.875
$ws.Range("H99").Value = 1775.125
$ws.Range("I99").Value = 1800.2858
$ws.Range("K99").Value = 1800.2858
$ws.Range("M99").Value = -302.2858000000001
$ws.Range("H126").Value = 1775.125
$ws.Range("I126").Value = 1800.2858
$ws.Range("K126").Value = 5400.857400000001
$ws.Range("M126").Value = -2930.857400000001
$ws.Range("H129").Value = 94500
$ws.Range("J129").Value = 94500
$ws.Range("L129").Value = 94500
$ws.Range("N129").Value = -104500
$ws.Range("H134").Value = 12989717
$ws.Range("I134").Value = 17859610
$ws.Range("J134").Value = 3332
$ws.Range("K134").Value = 53578830
$ws.Range("L134").Value = 9996
$ws.Range("M134").Value = -53576295
$ws.Range("N134").Value = -15066
$ws.Range("H136").Value = 2745.2222
$ws.Range("I136").Value = 2650.875
$ws.Range("K136").Value = 7952.625
$ws.Range("M136").Value = -5402.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1197212.5
$ws.Range("I4").Value = 646830
$ws.Range("J4").Value = 10003333
$ws.Range("K4").Value = 1940490
$ws.Range("L4").Value = 30009999
$ws.Range("M4").Value = -1940378
$ws.Range("N4").Value = -30010223
$ws.Range("H33").Value = 165
$ws.Range("J33").Value = 150
$ws.Range("L33").Value = 900
$ws.Range("N33").Value = -1466
$ws.Range("H34").Value = 3145.2632
$ws.Range("J34").Value = 3471.2354
$ws.Range("L34").Value = 10413.7062
$ws.Range("N34").Value = -10581.7062
$ws.Range("H55").Value = 1372.5
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H80").Value = 3329.3333
$ws.Range("J80").Value = 3344.5
$ws.Range("L80").Value = 10033.5
$ws.Range("N80").Value = -11905.5
$ws.Range("H83").Value = 3329.3333
$ws.Range("J83").Value = 3344.5
$ws.Range("L83").Value = 30100.5
$ws.Range("N83").Value = -39460.5
$ws.Range("H113").Value = 882.9167
$ws.Range("I113").Value = 876.25
$ws.Range("J113").Value = 886.25
$ws.Range("K113").Value = 2628.75
$ws.Range("L113").Value = 2658.75
$ws.Range("M113").Value = -458.75
$ws.Range("N113").Value = -6998.75
$ws.Range("H114").Value = 1028
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H136").Value = 5899.5
$ws.Range("J136").Value = 10000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -40200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 103599.37
$ws.Range("I14").Value = 157727.72
$ws.Range("J14").Value = 8874.75
$ws.Range("K14").Value = 157727.72
$ws.Range("L14").Value = 8874.75
$ws.Range("M14").Value = -157559.72
$ws.Range("N14").Value = -9210.75
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H80").Value = 1546.1666
$ws.Range("I80").Value = 1555.6
$ws.Range("J80").Value = 1499
$ws.Range("K80").Value = 1555.6
$ws.Range("L80").Value = 1499
$ws.Range("M80").Value = -557.5999999999999
$ws.Range("N80").Value = -3495
$ws.Range("H83").Value = 1546.1666
$ws.Range("I83").Value = 1555.6
$ws.Range("J83").Value = 1499
$ws.Range("K83").Value = 7778
$ws.Range("L83").Value = 7495
$ws.Range("M83").Value = -2786
$ws.Range("N83").Value = -17479
$ws.Range("H126").Value = 2787.0833
$ws.Range("J126").Value = 2975
$ws.Range("L126").Value = 8925
$ws.Range("N126").Value = -13865
$ws.Range("H132").Value = 2149.4517
$ws.Range("I132").Value = 1858.4615
$ws.Range("J132").Value = 3662.6
$ws.Range("K132").Value = 5575.3845
$ws.Range("L132").Value = 10987.8
$ws.Range("M132").Value = -3045.3845
$ws.Range("N132").Value = -16047.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2066.6924
$ws.Range("I7").Value = 2066.6924
$ws.Range("K7").Value = 2066.6924
$ws.Range("M7").Value = -1954.6924
$ws.Range("H22").Value = 1787
$ws.Range("I22").Value = 2499
$ws.Range("J22").Value = 1549.6666
$ws.Range("K22").Value = 2499
$ws.Range("L22").Value = 1549.6666
$ws.Range("M22").Value = -2204
$ws.Range("N22").Value = -2139.6666
$ws.Range("H27").Value = 1787
$ws.Range("I27").Value = 2499
$ws.Range("J27").Value = 1549.6666
$ws.Range("K27").Value = 2499
$ws.Range("L27").Value = 1549.6666
$ws.Range("M27").Value = -2392
$ws.Range("N27").Value = -1763.6666
$ws.Range("H40").Value = 3888.0667
$ws.Range("I40").Value = 3864.1853
$ws.Range("K40").Value = 3864.1853
$ws.Range("M40").Value = -3728.1853
$ws.Range("H126").Value = 2066.6924
$ws.Range("I126").Value = 2066.6924
$ws.Range("K126").Value = 6200.0772
$ws.Range("M126").Value = -3730.0772
$ws.Range("H128").Value = 80000
$ws.Range("J128").Value = 80000
$ws.Range("L128").Value = 80000
$ws.Range("N128").Value = -89960
$ws.Range("H130").Value = 120000
$ws.Range("J130").Value = 120000
$ws.Range("L130").Value = 120000
$ws.Range("N130").Value = -130040
$ws.Range("H132").Value = 3543.4814
$ws.Range("I132").Value = 4558.5
$ws.Range("J132").Value = 2946.4119
$ws.Range("K132").Value = 13675.5
$ws.Range("L132").Value = 8839.235700000001
$ws.Range("M132").Value = -11145.5
$ws.Range("N132").Value = -13899.2357

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H113").Value = 429.25925
$ws.Range("I113").Value = 405.21054
$ws.Range("J113").Value = 486.375
$ws.Range("K113").Value = 1215.63162
$ws.Range("L113").Value = 1459.125
$ws.Range("M113").Value = 954.3683800000001
$ws.Range("N113").Value = -5799.125
$ws.Range("H122").Value = 3397.75
$ws.Range("I122").Value = 3644
$ws.Range("K122").Value = 10932
$ws.Range("M122").Value = -8482
$ws.Range("H126").Value = 1692.9445
$ws.Range("I126").Value = 1631.25
$ws.Range("K126").Value = 4893.75
$ws.Range("M126").Value = -2423.75
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 2822.9656
$ws.Range("I132").Value = 2809.5
$ws.Range("K132").Value = 8428.5
$ws.Range("M132").Value = -5898.5
